$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right before the existing row 42 (old rows 42-70
# shift down to become rows 45-73).
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()
$ws.Rows.Item(42).Insert()

$newDate = Get-Date -Year 2021 -Month 12 -Day 29 -Hour 0 -Minute 0 -Second 0

# Row 42: Extra
$ws.Cells.Item(42,1).Value = 11
$ws.Cells.Item(42,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(42,3).Value = "Bíobío"
$ws.Cells.Item(42,4).Value = $newDate
$ws.Cells.Item(42,5).Value = 8
$ws.Cells.Item(42,6).Value = 100112028
$ws.Cells.Item(42,7).Value = "Sandia"
$ws.Cells.Item(42,8).Value = "Sin especificar"
$ws.Cells.Item(42,9).Value = "Extra"
$ws.Cells.Item(42,10).Value = 300
$ws.Cells.Item(42,11).Value = 2800
$ws.Cells.Item(42,12).Value = 2800
$ws.Cells.Item(42,13).Value = 2800
$ws.Cells.Item(42,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(42,15).Value = "Región de O'Higgins"
$ws.Cells.Item(42,16).Value = 2800
$ws.Cells.Item(42,17).Value = 1
$ws.Cells.Item(42,18).Value = "Hortaliza"

# Row 43: Primera
$ws.Cells.Item(43,1).Value = 11
$ws.Cells.Item(43,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(43,3).Value = "Bíobío"
$ws.Cells.Item(43,4).Value = $newDate
$ws.Cells.Item(43,5).Value = 8
$ws.Cells.Item(43,6).Value = 100112028
$ws.Cells.Item(43,7).Value = "Sandia"
$ws.Cells.Item(43,8).Value = "Sin especificar"
$ws.Cells.Item(43,9).Value = "Primera"
$ws.Cells.Item(43,10).Value = 300
$ws.Cells.Item(43,11).Value = 2400
$ws.Cells.Item(43,12).Value = 2400
$ws.Cells.Item(43,13).Value = 2400
$ws.Cells.Item(43,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(43,15).Value = "Región de O'Higgins"
$ws.Cells.Item(43,16).Value = 2400
$ws.Cells.Item(43,17).Value = 1
$ws.Cells.Item(43,18).Value = "Hortaliza"

# Row 44: Segunda
$ws.Cells.Item(44,1).Value = 11
$ws.Cells.Item(44,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(44,3).Value = "Bíobío"
$ws.Cells.Item(44,4).Value = $newDate
$ws.Cells.Item(44,5).Value = 8
$ws.Cells.Item(44,6).Value = 100112028
$ws.Cells.Item(44,7).Value = "Sandia"
$ws.Cells.Item(44,8).Value = "Sin especificar"
$ws.Cells.Item(44,9).Value = "Segunda"
$ws.Cells.Item(44,10).Value = 300
$ws.Cells.Item(44,11).Value = 2000
$ws.Cells.Item(44,12).Value = 2000
$ws.Cells.Item(44,13).Value = 2000
$ws.Cells.Item(44,14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(44,15).Value = "Región de O'Higgins"
$ws.Cells.Item(44,16).Value = 2000
$ws.Cells.Item(44,17).Value = 1
$ws.Cells.Item(44,18).Value = "Hortaliza"

Write-Output "Inserted rows 42-44; dimension now $($ws.UsedRange.Address())"
